# Generate Report for Handoff
#
# Inserts two new "Ready for handoff" file entries - 03d65910-50c3-439a-9026-2e189cbc29ee
# and 980bbc2d-6578-42db-8dd1-c2d8732bde02 - ahead of the existing last entry
# (b937ea7f-d83a-4a0a-9424-25d276c7f8b2) on all three sheets (Overview, zh-cn, de-de).
# The previously-last row's data is preserved and simply pushed down to the new
# bottom row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 5 now describes 03d65910-...
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/03d65910/e2e/03d65910-50c3-439a-9026-2e189cbc29ee.md", "", "", "03d65910-50c3-439a-9026-2e189cbc29ee.md") | Out-Null
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"
$ws1.Range("D5").Value2 = "2016-03-24 08:16:48"
$ws1.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 6 (new) describes 980bbc2d-...
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/980bbc2d/e2e/980bbc2d-6578-42db-8dd1-c2d8732bde02.md", "", "", "980bbc2d-6578-42db-8dd1-c2d8732bde02.md") | Out-Null
$ws1.Range("B6").Value2 = "Ready for handoff"
$ws1.Range("C6").Value2 = "Ready for handoff"
$ws1.Range("D6").Value2 = "2016-03-24 08:16:48"
$ws1.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 7 (new) is the original row-5 content (b937ea7f-...), moved down
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8d90920953fea1d8bf08c3695b91db9017de9b86/e2e/b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md", "", "", "b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md") | Out-Null
$ws1.Range("B7").Value2 = "Ready for handoff"
$ws1.Range("C7").Value2 = "Ready for handoff"
$ws1.Range("D7").Value2 = "2016-03-24 08:13:37"
$ws1.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 5 now describes 03d65910-...
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/03d65910/e2e/03d65910-50c3-439a-9026-2e189cbc29ee.md", "", "", "03d65910-50c3-439a-9026-2e189cbc29ee.md") | Out-Null
$ws2.Range("B5").Value2 = ".md"
$ws2.Range("C5").Value2 = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03d65910/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03d65910-50c3-439a-9026-2e189cbc29ee.5fc483cbee2d00524b28de3e83772b6453eea09b.zh-cn.xlf", "", "", "03d65910-50c3-439a-9026-2e189cbc29ee.5fc483cbee2d00524b28de3e83772b6453eea09b.zh-cn.xlf") | Out-Null
$ws2.Range("E5").Value2 = "2016-03-24 08:16:44"
$ws2.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J5").Value2 = "Include"

# Row 6 (new) describes 980bbc2d-...
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/980bbc2d/e2e/980bbc2d-6578-42db-8dd1-c2d8732bde02.md", "", "", "980bbc2d-6578-42db-8dd1-c2d8732bde02.md") | Out-Null
$ws2.Range("B6").Value2 = ".md"
$ws2.Range("C6").Value2 = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/980bbc2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/980bbc2d-6578-42db-8dd1-c2d8732bde02.49a22308174e0cfc8474abb63da51bb321dcc8fe.zh-cn.xlf", "", "", "980bbc2d-6578-42db-8dd1-c2d8732bde02.49a22308174e0cfc8474abb63da51bb321dcc8fe.zh-cn.xlf") | Out-Null
$ws2.Range("E6").Value2 = "2016-03-24 08:16:44"
$ws2.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H6").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J6").Value2 = "Include"

# Row 7 (new) is the original row-5 content (b937ea7f-...), moved down
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8d90920953fea1d8bf08c3695b91db9017de9b86/e2e/b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md", "", "", "b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md") | Out-Null
$ws2.Range("B7").Value2 = ".md"
$ws2.Range("C7").Value2 = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b47d768d822f5612d4e793b70b840302b4cdef5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b937ea7f-d83a-4a0a-9424-25d276c7f8b2.baab81d2c15e940567b209d48f888f7a38ce624b.zh-cn.xlf", "", "", "b937ea7f-d83a-4a0a-9424-25d276c7f8b2.baab81d2c15e940567b209d48f888f7a38ce624b.zh-cn.xlf") | Out-Null
$ws2.Range("E7").Value2 = "2016-03-24 08:13:33"
$ws2.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H7").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J7").Value2 = "Include"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 5 now describes 03d65910-...
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/03d65910/e2e/03d65910-50c3-439a-9026-2e189cbc29ee.md", "", "", "03d65910-50c3-439a-9026-2e189cbc29ee.md") | Out-Null
$ws3.Range("B5").Value2 = ".md"
$ws3.Range("C5").Value2 = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/03d65910/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03d65910-50c3-439a-9026-2e189cbc29ee.5fc483cbee2d00524b28de3e83772b6453eea09b.de-de.xlf", "", "", "03d65910-50c3-439a-9026-2e189cbc29ee.5fc483cbee2d00524b28de3e83772b6453eea09b.de-de.xlf") | Out-Null
$ws3.Range("E5").Value2 = "2016-03-24 08:16:48"
$ws3.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H5").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J5").Value2 = "Include"

# Row 6 (new) describes 980bbc2d-...
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/980bbc2d/e2e/980bbc2d-6578-42db-8dd1-c2d8732bde02.md", "", "", "980bbc2d-6578-42db-8dd1-c2d8732bde02.md") | Out-Null
$ws3.Range("B6").Value2 = ".md"
$ws3.Range("C6").Value2 = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/980bbc2d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/980bbc2d-6578-42db-8dd1-c2d8732bde02.49a22308174e0cfc8474abb63da51bb321dcc8fe.de-de.xlf", "", "", "980bbc2d-6578-42db-8dd1-c2d8732bde02.49a22308174e0cfc8474abb63da51bb321dcc8fe.de-de.xlf") | Out-Null
$ws3.Range("E6").Value2 = "2016-03-24 08:16:48"
$ws3.Range("E6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H6").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J6").Value2 = "Include"

# Row 7 (new) is the original row-5 content (b937ea7f-...), moved down
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8d90920953fea1d8bf08c3695b91db9017de9b86/e2e/b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md", "", "", "b937ea7f-d83a-4a0a-9424-25d276c7f8b2.md") | Out-Null
$ws3.Range("B7").Value2 = ".md"
$ws3.Range("C7").Value2 = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b62b93eda970aa064249d4a6373ffcf5fbbae491/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b937ea7f-d83a-4a0a-9424-25d276c7f8b2.baab81d2c15e940567b209d48f888f7a38ce624b.de-de.xlf", "", "", "b937ea7f-d83a-4a0a-9424-25d276c7f8b2.baab81d2c15e940567b209d48f888f7a38ce624b.de-de.xlf") | Out-Null
$ws3.Range("E7").Value2 = "2016-03-24 08:13:37"
$ws3.Range("E7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H7").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J7").Value2 = "Include"

Write-Output "Generated handoff report rows for 03d65910-50c3-439a-9026-2e189cbc29ee and 980bbc2d-6578-42db-8dd1-c2d8732bde02 across Overview, zh-cn, de-de."
